# Add a new "債務" (debt) worksheet after the existing "存款" (deposit) sheet,
# carrying the same header/data-row formatting as sheet1, then fill in the
# new header row and single data row describing a 農會貸款 (credit union loan)
# debt entry.

$xlPasteFormats = -4122
$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after "存款" so tab order becomes 存款, 債務.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "債務"

# --- Header row (row 1): reuse sheet1's header style (bold + border) ---
# sheet1 only spans columns A:M (13 cols); sheet2 needs a 14th ("index")
# column, so N1 borrows sheet1's M1 ("index" header) formatting as its
# style source.
$headers = [ordered]@{
    "B1" = @{ text = "species";            styleFrom = "B1" };
    "C1" = @{ text = "debtor";             styleFrom = "C1" };
    "D1" = @{ text = "owner";              styleFrom = "D1" };
    "E1" = @{ text = "total";              styleFrom = "E1" };
    "F1" = @{ text = "register_date";      styleFrom = "F1" };
    "G1" = @{ text = "register_reason";    styleFrom = "G1" };
    "H1" = @{ text = "property_category";  styleFrom = "H1" };
    "I1" = @{ text = "category";           styleFrom = "I1" };
    "J1" = @{ text = "date";               styleFrom = "J1" };
    "K1" = @{ text = "legislator_name";    styleFrom = "K1" };
    "L1" = @{ text = "legislator_id";      styleFrom = "L1" };
    "M1" = @{ text = "source_file";        styleFrom = "M1" };
    "N1" = @{ text = "index";              styleFrom = "M1" };
}

foreach ($addr in $headers.Keys) {
    $info = $headers[$addr]
    $ws1.Range($info.styleFrom).Copy()
    $ws2.Range($addr).PasteSpecial($xlPasteFormats)
    $ws2.Range($addr).Value = $info.text
}

# --- Data row (row 2): reuse sheet1's data style (N2 borrows sheet1's M2) ---
$rowTwoText = [ordered]@{
    "B2" = @{ text = "農會貸款";        styleFrom = "B2" };
    "C2" = @{ text = "劉建國";          styleFrom = "C2" };
    "D2" = @{ text = "雲林縣斗六市農會雲林縣斗六市民生路"; styleFrom = "D2" };
    "G2" = @{ text = "代償債務";        styleFrom = "G2" };
    "H2" = @{ text = "debt";            styleFrom = "H2" };
    "I2" = @{ text = "normal";          styleFrom = "I2" };
    "K2" = @{ text = "劉建國";          styleFrom = "K2" };
    "M2" = @{ text = "tmpd6c01";        styleFrom = "M2" };
}

foreach ($addr in $rowTwoText.Keys) {
    $info = $rowTwoText[$addr]
    $ws1.Range($info.styleFrom).Copy()
    $ws2.Range($addr).PasteSpecial($xlPasteFormats)
    $ws2.Range($addr).Value = $info.text
}

# F2 holds a "99年12月08日"-style ROC date string; it isn't auto-parsed as a
# date by the host, so a plain format-copy + Value assignment is fine.
$ws1.Range("F2").Copy()
$ws2.Range("F2").PasteSpecial($xlPasteFormats)
$ws2.Range("F2").Value = "99年12月08日"

# J2 holds "2012-05-01" which Excel *would* auto-parse into a serial date if
# assigned through .Value (it matches yyyy-mm-dd). sheet1's I2 cell already
# stores this exact string as literal text, so paste its typed value
# straight across (xlPasteValues) after copying J2's own formatting — this
# keeps the cell text-typed without inventing a new text number format.
$ws1.Range("J2").Copy()
$ws2.Range("J2").PasteSpecial($xlPasteFormats)
$ws1.Range("I2").Copy()
$ws2.Range("J2").PasteSpecial($xlPasteValues)

$rowTwoNumber = [ordered]@{
    "A2" = @{ value = 90;      styleFrom = "A2" };
    "E2" = @{ value = 5200000; styleFrom = "E2" };
    "L2" = @{ value = 1723;    styleFrom = "L2" };
    "N2" = @{ value = 90;      styleFrom = "M2" };
}

foreach ($addr in $rowTwoNumber.Keys) {
    $info = $rowTwoNumber[$addr]
    $ws1.Range($info.styleFrom).Copy()
    $ws2.Range($addr).PasteSpecial($xlPasteFormats)
    $ws2.Range($addr).Value = $info.value
}

$excel.CutCopyMode = $false

# Keep the first sheet active/selected, matching the original workbook.
$ws1.Select()
